$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (account for the ColumnWidth <-> stored-width offset)
$ws.Columns.Item(1).ColumnWidth = 5.166666666666667
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668

# Update row 2 values (keep header row 1 and E2 unchanged)
$ws.Range("A2").Value = "Você"
$ws.Range("B2").Value = "Desconhecido"
$ws.Range("C2").Value = "R$ 100,00"
$ws.Range("D2").Value = "JOSIE ADAUANE DIAS"

# Remove rows 3 to 6 entirely
$ws.Range("A3:E6").Delete()
